$d = $word.ActiveDocument

# --- Step 1: fix "I is" -> "It is" (the grammar-check markers around "is" no longer apply) ---
$null = $d.Content.Find.Execute(
    " I is a Clydesdale mare named Millie.", $true, $false, $false, $false, $false,
    $true, 1, $false, " It is a Clydesdale mare named Millie.", 2)

# --- Locate the "Millie" paragraph; this is where the new scene gets inserted ---
$millie = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*named Millie.*") {
        $millie = $i
        break
    }
}

# --- Step 2: drop 9 of the 11 blank paragraphs that used to trail the chapter (2 remain at the end) ---
$toDelete = $d.Paragraphs.Count - $millie - 2
for ($n = 0; $n -lt $toDelete; $n++) {
    $d.Paragraphs($millie + 1).Range.Delete()
}

# --- Step 3: insert the new paragraphs of story content (end of Ch.2 + Chapter 3 opening) ---
$texts = @(
    "“You haven’t said much.” Pony says.",
    "“Yeah, I know. I have a bit on my mind.”",
    "“Like what?”",
    "I pull out the ring and show her. “Somebody named Lela gave me this.”",
    "“Did you say Lela?”",
    "“Yeah why?”. I asked.",
    "She throws the saddle cloth on and boosts me up. “We’re going to my house.” Pony says.",
    "“Alright?” is something wrong?",
    "“Yeah.” I say showing her the ring.",
    "“I know someone named Lela that had a similar ring. She was the best friend a friend could have, but she died about a year ago. “She leads the horse.",
    "“What happened to her?”",
    "Pony hesitates then says. “You wouldn’t understand. Did she give you a copper chain?” I untuck it from my shirt to show her, “Good, you wanna keep that on.”",
    "“Why?” I ask feeling inquisitive.",
    "“It protects you from being detected by magic.”",
    "“You went crazy too?”",
    "“No. I will prove it when we get home.”",
    "For the rest of the ride there was an uneasy silence between us.",
    "It takes us about 20 minutes for us to get through the front gate of her small ranch. Toi the right is a 10 acre and on the left there’s her two-story house. She puts her horse in the pasture.",
    "“Let me see the ring.” I hand it to her. She takes the ring and lays it on her palm then places her other hand on it and says “be seen” then hold it up between her forefinger and thumb. I see the other would that I had been in before. I can’t see it well enough to make anything out because she is holding it. She puts the ring between her hands as before and says, “be blind”. Now it is the way it was. She goes inside without saying anything and comes back out a couple minutes later with an old-fashioned key, it looks new though.",
    "“But is that where Lela is?” I finally ask to break the silence.",
    "“It should be. I hope it is. I really do!” She said in the most ecstatic voice I have ever hear. Shye continues, “: I can’t wait to see her again, she’s really nice. She gave me her horse!” she exclaims. “How do I get there?” I ask.",
    "“With this key.” She says turning toward her house door. She unlocks her door with the key and opens it, the inside of the house is trashed wood floors are rotting, most windows are broke and the ones that aren’t are moldy. The table and chairs are scattered everywhere. She takes a step in and motions for me to follow. I cautiously take a step in and she shuts the door and locks it with the key, bouts the key in her pocket and unlatches is with her hand and reopens the door.",
    "",
    "Chapter 3",
    "The grass outside is all dead and in the field is mostly burned, some trees are still smoking from the fire. Ash is floating in the air. The porch had been broken and has holed in the top. The white fence that kept the horses in was broken and scattered everywhere. While they were taking in all the information of this world a large gray reptilian demon creature crawls out from beside the house. It is about three feet in length, has bat-like wings and sharp teeth, and spines all over its body, standing up right. It runs across the yard to the girls and tries to bite at Pony, Taylor screams. Pony gets mad at it and stomps its head with her boot. It squirms in pain; Taylor is covering her eyes with her hands. The creature gets us quickly and runs off. Pony starts off walking into the field ash blowing up as she sleeps. I stand there for a moment in fear before realizing she is leaving me behind. I gather my courage in hopes that it is a one-time event and call out “Hey wait up!” and I make a small dash up next to her.",
    "I look to her in disbelief, “Where are we.”",
    "“Lela calls it the Gray Area. There’s not supposed to be any life here, something must have gone wrong when we entered. Hopefully there is nothing else here, if there is it has to be removed.”",
    "“Where exactly are we going?” I ask.",
    "“We have to go back to the school.”",
    "“But we don’t usually have to have to go through the field to get to the school. So why are we going this way?” I ask inquisitively.",
    "“To many questions will ruin anything.” She replies.",
    "I stop asking questions and start observing the ruined landscape. After awhile I begin to realize that there really are no animals of any kind here. The air is thick with ash now from us walking in the field. A few moments later Pony’s boot strikes a single plank with a keyhole, she takes her old key and “unlocks” the plank. A door reveals a passage leading to a small room underground. The fold already has a lantern inside, but I can’t tell what’s in there. She hops in, I just wait on the outside, occasionally looking around to see if anything was around. Clank. Pony tosses out an old rusted steel sword. The sudden notice spooked me, now my heart is pounding. Pony climbs out of the hold "
)
$borderFlags = @(
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false
)

$anchorIdx = $millie
for ($i = 0; $i -lt $texts.Length; $i++) {
    $d.Paragraphs($anchorIdx).Range.InsertParagraphAfter()
    $anchorIdx = $anchorIdx + 1
    $newPara = $d.Paragraphs($anchorIdx)
    if ($texts[$i].Length -gt 0) {
        $newPara.Range.InsertAfter($texts[$i])
    }
    if ($borderFlags[$i]) {
        $newPara.Borders(-3).LineStyle = 1
        $newPara.Borders(-3).LineWidth = 3
        $newPara.Borders.DistanceFromBottom = 1
        $newPara.Borders(-3).ColorIndex = 0
    }
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)